# Update the player roster table (A2:C19) to the new data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Chris Paul",        "PG",       "San Antonio Spurs"),
    @("Stephon Castle",    "PG,SG",    "San Antonio Spurs"),
    @("Jaylen Brown",      "SG,SF",    "Boston Celtics"),
    @("Payton Pritchard",  "PG,SG",    "Boston Celtics"),
    @("Deni Avdija",       "SF,PF",    "Portland Trail Blazers"),
    @("Chet Holmgren",     "PF,C",     "Oklahoma City Thunder"),
    @("Pascal Siakam",     "SF,PF,C",  "Indiana Pacers"),
    @("Dillon Brooks",     "SG,SF,PF", "Houston Rockets"),
    @("Nikola Jokic",      "C",        "Denver Nuggets"),
    @("Nick Richards",     "C",        "Phoenix Suns"),
    @("Cole Anthony",      "PG",       "Orlando Magic"),
    @("Rudy Gobert",       "C",        "Minnesota Timberwolves"),
    @("Jalen Green",       "PG,SG",    "Houston Rockets"),
    @("Paolo Banchero",    "SF,PF",    "Orlando Magic"),
    @("Ayo Dosunmu",       "PG,SG,SF", "Chicago Bulls"),
    @("Russell Westbrook", "PG,SG",    "Denver Nuggets"),
    @("Jalen Suggs",       "PG,SG",    "Orlando Magic"),
    @("Jakob Poeltl",      "C",        "Toronto Raptors")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
